$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (Date serial, Method, ElapsedMs, wordCount, sentenceCount,
# posWordCount, negWordCount, posWordPercentage, negWordPercentage,
# positivePhraseCount, negativePhraseCount, posPhrasePercentage, negPhrasePercentage)
$rows = @(
    @(42600.792222222219, "Named", 11630, 5242, 270, 63, 19, 76, 22, 2, 0, 100, 0),
    @(42600.794571759259, "Named", 15131, 5925, 327, 63, 25, 71, 28, 2, 0, 100, 0),
    @(42600.830810185187, "Named",  8898, 4543, 238, 46, 23, 66, 33, 2, 0, 100, 0),
    @(42600.879212962966, "Named",  7882, 3767, 192, 30, 23, 56, 43, 1, 0, 100, 0)
)

$startRow = 3
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $data[0]
    $aCell.NumberFormat = "m/d/yy h:mm"

    for ($c = 1; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}

# Widen column A slightly to fit the new, longer date/time content
$ws.Columns.Item(1).ColumnWidth = 14
